# The edit reorders the data rows (rows 2-119, i.e. the weekly blocks of
# price records) into a new sequence while keeping every record's own
# values (date, quality, volumes, prices, origin, etc.) fully intact.
# No cell values besides their row position actually change.
#
# Build the old-row -> new-row mapping (1-based spreadsheet rows) and use
# a single bulk read + bulk write through Value2 so every column (A-R)
# moves together as a unit for each record.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstDataRow = 2
$lastDataRow = 119
$lastCol = 18   # column R

# For each new row (index into this array, offset by $firstDataRow),
# this holds which OLD row's data should be placed there.
$sourceRowForNewRow = @(101,102,103,104,105,54,55,56,47,48,68,69,70,71,115,116,117,118,63,64,65,100,96,97,12,106,107,108,32,44,91,92,93,11,45,46,66,52,53,77,26,27,28,82,83,2,3,4,5,6,7,8,38,39,40,67,29,30,31,95,88,89,90,78,79,80,81,85,86,87,41,42,43,49,50,51,25,10,16,75,76,109,110,111,112,119,20,17,18,19,21,22,23,24,33,34,35,36,37,59,60,61,9,13,14,15,98,99,57,58,62,94,84,113,114,72,73,74)

# Read the entire original block once (values only, Value2 keeps Excel
# date serials as numbers so the D-column style still renders them as
# dates afterwards).
$sourceRange = $ws.Range("A$($firstDataRow):R$lastDataRow")
$original = $sourceRange.Value2

$rowCount = $lastDataRow - $firstDataRow + 1

# Build the reordered 2D array.
$reordered = New-Object 'object[,]' $rowCount, $lastCol
for ($i = 0; $i -lt $rowCount; $i++) {
    $srcRow = $sourceRowForNewRow[$i]
    $srcIndex = $srcRow - $firstDataRow + 1   # 1-based index into $original
    for ($c = 1; $c -le $lastCol; $c++) {
        $reordered[$i, $c - 1] = $original[$srcIndex, $c]
    }
}

# Write the reordered block back in one shot.
$destRange = $ws.Range("A$($firstDataRow):R$lastDataRow")
$destRange.Value2 = $reordered
